# Updates the "Rider_A" transition-probability matrix (Sheet1) to reflect
# additional simulated games: row totals (and therefore each row's
# count/total probabilities) were recomputed after more games were run.
# Rows 2-19 correspond to states Af0..Br0 (column A), columns B-S are the
# 18 target states. Row 1 (headers) and column A (state labels) are
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(0.2007722007722008,0.5752895752895753,0,0,0,0,0,0,0.01158301158301158,0,0,0,0,0,0.1351351351351351,0,0,0.07722007722007722),
    @(0.006535947712418301,0.0196078431372549,0,0,0,0,0,0,0.03267973856209151,0,0,0,0,0,0.7777777777777778,0,0,0.1633986928104575),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.7837837837837838,0,0,0.2162162162162162),
    @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.6666666666666666,0,0,0.3333333333333333),
    @(0.04411764705882353,0,0.01470588235294118,0.004901960784313725,0.04901960784313725,0,0,0,0.2892156862745098,0,0,0,0,0.04411764705882353,0,0.1568627450980392,0.06372549019607843,0.3333333333333333),
    @(0.07567567567567568,0,0.02702702702702703,0,0.04324324324324325,0,0,0,0.1567567567567568,0,0,0,0,0.01621621621621622,0,0.1945945945945946,0.07567567567567568,0.4108108108108108),
    @(0.078125,0,0.01822916666666667,0.002604166666666667,0.05989583333333334,0,0,0,0.125,0,0,0,0,0.005208333333333333,0,0.2109375,0.07291666666666667,0.4270833333333333),
    @(0.08762886597938144,0,0.02061855670103093,0,0.08762886597938144,0,0,0,0.09278350515463918,0,0,0,0,0.01030927835051546,0,0.2061855670103093,0.06185567010309279,0.4329896907216495),
    @(0.1112033195020747,0,0.01825726141078838,0.0008298755186721991,0.06721991701244813,0,0,0,0.1286307053941909,0,0,0,0,0.01410788381742739,0,0.2298755186721992,0.06639004149377593,0.3634854771784232),
    @(0,0,0,0,0,0.1482758620689655,0,0,0.1068965517241379,0.1931034482758621,0.5344827586206896,0,0,0,0,0,0,0.01724137931034483),
    @(0,0,0,0,0,0.7222222222222222,0,0,0.2222222222222222,0.01234567901234568,0.01851851851851852,0,0,0,0,0,0,0.02469135802469136),
    @(0,0,0,0,0,0.6976744186046512,0,0,0.2558139534883721,0,0,0,0,0,0,0,0,0.04651162790697674),
    @(0,0,0,0,0,1.0,0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0.02094240837696335,0,0.1361256544502618,0.08376963350785341,0.3507853403141361,0.06282722513089005,0,0.01047120418848168,0,0.08900523560209424,0,0,0,0.2460732984293194),
    @(0,0,0,0,0.01648351648351648,0,0.1758241758241758,0.08241758241758242,0.4615384615384616,0.1098901098901099,0,0.01098901098901099,0,0.06593406593406594,0,0,0,0.07692307692307693),
    @(0,0,0,0,0.01079913606911447,0,0.1490280777537797,0.1036717062634989,0.4514038876889849,0.07775377969762419,0,0.02591792656587473,0,0.06047516198704104,0,0,0,0.1209503239740821),
    @(0,0,0,0,0.03472222222222222,0,0.2430555555555556,0.1041666666666667,0.3958333333333333,0.09722222222222222,0,0.02083333333333333,0,0.04166666666666666,0,0,0,0.0625),
    @(0,0,0,0,0.01913043478260869,0,0.1973913043478261,0.08695652173913043,0.3608695652173913,0.1252173913043478,0,0.02086956521739131,0.0008695652173913044,0.06347826086956522,0,0,0,0.1252173913043478)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $rowData = $newValues[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $rowData[$j]
    }
}
